# Auto-update GitHub repos Excel export
# The "Project" repo entry (previously row 12) was removed from the list,
# shifting the remaining rows (repo-scanner, skill_captain_dsa_arrays,
# The-Warrior) up by one and shrinking the used range by a row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 12 ("Project" / its URL) and shift the rows below it
# (13, 14, 15) up by one.
$ws.Rows(12).Delete()
